$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimpleSearch")

# Fix the typo in the variable name: "Reports__Tab_URL" -> "Reports_Tab_URL"
$ws.Range("A10").Value = "Reports_Tab_URL"

# Update the active cell selection on the sheet to A9
$ws.Activate()
$ws.Range("A9").Select()
